# Add more error messages
# - Replace the "msgBadAttachments" error entry on the ExceptionMessages sheet
#   with a new "msgBadAttachment" entry that has an updated message.
# - Make the ExceptionMessages sheet the active sheet (with a new selection),
#   instead of the Settings sheet.

$wb = $excel.ActiveWorkbook

$wsExceptions = $wb.Worksheets.Item("ExceptionMessages")

# Update the (previously) "msgBadAttachments" row with the new key/message.
$wsExceptions.Range("A3").Value = "msgBadAttachment"
$wsExceptions.Range("B3").Value = "Email does not contain exactly one xlsx file"

# Make ExceptionMessages the active/selected sheet, with D11 selected.
$wsExceptions.Activate()
$wsExceptions.Range("D11").Select()
